$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 53453
$ws.Range("B2").Value = "Enzo Novaes"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45100
$ws.Range("G2").Value = 6619.01

$ws.Range("A3").Value = 10669
$ws.Range("B3").Value = "Ian da Paz"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45100
$ws.Range("G3").Value = 6177.58

$ws.Range("A4").Value = 78241
$ws.Range("B4").Value = "Carlos Eduardo da Luz"
$ws.Range("C4").Value = "Operações"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45084
$ws.Range("G4").Value = 4514.14

$ws.Range("A5").Value = 77456
$ws.Range("B5").Value = "Sr. Vinicius Santos"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45091
$ws.Range("G5").Value = 11699.21

$ws.Range("A6").Value = 12050
$ws.Range("B6").Value = "Isaac Melo"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45090
$ws.Range("G6").Value = 9688.51

$ws.Range("A7").Value = 82903
$ws.Range("B7").Value = "Dra. Heloísa das Neves"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 10614.87

$ws.Range("A8").Value = 91864
$ws.Range("B8").Value = "Melissa da Mota"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45078
$ws.Range("G8").Value = 12364.49

$ws.Range("A9").Value = 47824
$ws.Range("B9").Value = "Augusto Pires"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 5686.56

$ws.Range("A10").Value = 35033
$ws.Range("B10").Value = "Helena Ferreira"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45096
$ws.Range("G10").Value = 3028.63

$ws.Range("A11").Value = 51478
$ws.Range("B11").Value = "Lucas Oliveira"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45097
$ws.Range("G11").Value = 10104.63
